# Weekly data refresh: insert two new daily price records at the top of the
# historical data block (rows 533-534), pushing the existing rows down by two
# (old row 533 -> new row 535, ... old row 591 -> new row 593).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 533; formats (incl. the
# date style used by column D) are copied down from row 532/534 automatically.
$ws.Rows("533:534").Insert()

# --- New row 533 -----------------------------------------------------------
$ws.Cells.Item(533, 1).Value  = 9
$ws.Cells.Item(533, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(533, 3).Value  = "Metropolitana"
$ws.Cells.Item(533, 4).Value  = 45106
$ws.Cells.Item(533, 5).Value  = 13
$ws.Cells.Item(533, 6).Value  = 100112032
$ws.Cells.Item(533, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(533, 8).Value  = "Bola 8"
$ws.Cells.Item(533, 9).Value  = "Primera"
$ws.Cells.Item(533, 10).Value = 52
$ws.Cells.Item(533, 11).Value = 12000
$ws.Cells.Item(533, 12).Value = 14000
$ws.Cells.Item(533, 13).Value = 13000
$ws.Cells.Item(533, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(533, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(533, 16).Value = 260
$ws.Cells.Item(533, 17).Value = 50
$ws.Cells.Item(533, 18).Value = "Hortaliza"

# --- New row 534 -----------------------------------------------------------
$ws.Cells.Item(534, 1).Value  = 9
$ws.Cells.Item(534, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(534, 3).Value  = "Metropolitana"
$ws.Cells.Item(534, 4).Value  = 45106
$ws.Cells.Item(534, 5).Value  = 13
$ws.Cells.Item(534, 6).Value  = 100112032
$ws.Cells.Item(534, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(534, 8).Value  = "Sin especificar"
$ws.Cells.Item(534, 9).Value  = "Primera"
$ws.Cells.Item(534, 10).Value = 70
$ws.Cells.Item(534, 11).Value = 12000
$ws.Cells.Item(534, 12).Value = 14000
$ws.Cells.Item(534, 13).Value = 13000
$ws.Cells.Item(534, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(534, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(534, 16).Value = 260
$ws.Cells.Item(534, 17).Value = 50
$ws.Cells.Item(534, 18).Value = "Hortaliza"
